# edit.ps1 - reproduce the tracked-edit diff via Word COM-interop
#
# Summary of changes:
#  1. Move the "_GoBack" bookmark from the "Single LEDs." paragraph to the
#     end of the "Introduction" heading paragraph.
#  2/7. Split "Here are some potential {input,output} devices ... )" into
#     3 runs with proofErr gramStart/gramEnd wrapping "( 1".
#  3/4/5. Split the PMOD joystick / orientation sensor / audio ADC table
#     cells so "Digilent" is wrapped in proofErr spellStart/spellEnd.
#  6/9. Split "... memory mapped i/o or a port based approach)? " runs so
#     "i/o" is wrapped in proofErr gramStart/gramEnd.
#  10. Split "... kind of like cin)" so "cin" is wrapped in proofErr
#     spellStart/spellEnd and the closing ")" becomes its own run.
#  11. Footer page-number field cached result "2" -> "1".

$d = $word.ActiveDocument

function Replace-RunXml {
    param(
        $Doc,
        [string]$SearchText,
        [string]$NewRunsXml,
        $StartFrom = $null
    )

    $seeker = $Doc.Content
    if ($null -ne $StartFrom) {
        $seeker.Start = $StartFrom
    }
    $found = $seeker.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $SearchText"
        return $null
    }

    $target = $Doc.Range($seeker.Start, $seeker.End)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $NewRunsXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
    return $target.Start
}

# ---------------------------------------------------------------------
# Edit 1: move the "_GoBack" bookmark
# ---------------------------------------------------------------------

# Remove it from its old spot (right after the "S" run in "Single LEDs.")
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-add it right after "...troduction" (end of the Introduction heading),
# but *before* the paragraph mark. A bookmark collapsed exactly on a
# paragraph-content boundary needs a temporary anchor character, or the
# engine mis-places it at the start of the document/paragraph.
$introPara = $d.Paragraphs(2)
$introRange = $introPara.Range
$endPos = $introRange.End - 1   # position right before the paragraph mark

$anchor = $d.Range($endPos, $endPos)
$anchor.InsertAfter("X")
$bmPos = $endPos
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Text = ""

# ---------------------------------------------------------------------
# Edit 2: "Here are some potential input devices ..." -> proofErr split
# ---------------------------------------------------------------------

$inputDevicesNew = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr>' +
    '<w:t xml:space="preserve">Here are some potential input devices and relative levels of difficulty </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>( 1</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr>' +
    '<w:t xml:space="preserve"> – easy , 3 – medium, 5 – hardest )</w:t></w:r>'

Replace-RunXml $d "Here are some potential input devices and relative levels of difficulty ( 1 – easy , 3 – medium, 5 – hardest )" $inputDevicesNew | Out-Null

# ---------------------------------------------------------------------
# Edit 3: PMOD joystick cell -> "Digilent" spellStart/spellEnd
# ---------------------------------------------------------------------

$pmodNew = '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr>' +
    '<w:t xml:space="preserve">The PMOD joysticks. These require talking over a serial interface (similar to what the keyboard does).  There is example code on </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr><w:t>Digilent</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'

Replace-RunXml $d "The PMOD joysticks. These require talking over a serial interface (similar to what the keyboard does).  There is example code on Digilent" $pmodNew | Out-Null

# ---------------------------------------------------------------------
# Edit 4: orientation sensors cell -> "Digilent" spellStart/spellEnd
# ---------------------------------------------------------------------

$orientationNew = '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr>' +
    '<w:t xml:space="preserve">The built in orientation sensors. There is example code on </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr><w:t>Digilent</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr><w:t>.</w:t></w:r>'

Replace-RunXml $d "The built in orientation sensors. There is example code on Digilent." $orientationNew | Out-Null

# ---------------------------------------------------------------------
# Edit 5: Audio ADC cell -> "Digilent" spellStart/spellEnd
# ---------------------------------------------------------------------

$audioNew = '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr>' +
    '<w:t xml:space="preserve">The Audio input Analog to Digital converter. There is example code on </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr><w:t>Digilent</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="18"/></w:rPr>' +
    '<w:t xml:space="preserve"> but you will need to spend effort to implement it in VHDL.</w:t></w:r>'

Replace-RunXml $d "The Audio input Analog to Digital converter. There is example code on Digilent but you will need to spend effort to implement it in VHDL." $audioNew | Out-Null

# ---------------------------------------------------------------------
# Edit 6: "How will you interface these input devices ..." -> "i/o" split
# ---------------------------------------------------------------------

$ioInputNew = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr>' +
    '<w:t xml:space="preserve">How will you interface these input devices to your processor (e.g. memory mapped </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr><w:t>i/o</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr>' +
    '<w:t xml:space="preserve"> or a port based approach)? </w:t></w:r>'

Replace-RunXml $d "How will you interface these input devices to your processor (e.g. memory mapped i/o or a port based approach)? " $ioInputNew | Out-Null

# ---------------------------------------------------------------------
# Edit 7: "Here are some potential output devices ..." -> proofErr split
# ---------------------------------------------------------------------

$outputDevicesNew = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr>' +
    '<w:t xml:space="preserve">Here are some potential output devices and relative levels of difficulty </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>( 1</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr>' +
    '<w:t xml:space="preserve"> – easy , 3 – medium, 5 – hardest )</w:t></w:r>'

Replace-RunXml $d "Here are some potential output devices and relative levels of difficulty ( 1 – easy , 3 – medium, 5 – hardest )" $outputDevicesNew | Out-Null

# ---------------------------------------------------------------------
# Edit 9: " devices to your processor (e.g. memory mapped i/o ...)? " (2nd
# occurrence -- the "output" devices paragraph) -> "i/o" split
# ---------------------------------------------------------------------

$ioOutputNew = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr>' +
    '<w:t xml:space="preserve"> devices to your processor (e.g. memory mapped </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr><w:t>i/o</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr>' +
    '<w:t xml:space="preserve"> or a port based approach)? </w:t></w:r>'

# Skip past the first (input-devices) occurrence by starting the search
# right after it.
$firstIoOutput = $d.Content
$firstIoOutput.Find.Execute("How will you interface these input devices", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Replace-RunXml $d " devices to your processor (e.g. memory mapped i/o or a port based approach)? " $ioOutputNew $firstIoOutput.End | Out-Null

# ---------------------------------------------------------------------
# Edit 10: " ... kind of like cin)" -> "cin" spellStart/spellEnd
# ---------------------------------------------------------------------

$cinNew = '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b w:val="0"/><w:color w:val="auto"/></w:rPr>' +
    '<w:t xml:space="preserve"> that allows your CPU to read text/numbers from a keyboard buffer(kind of like </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b w:val="0"/><w:color w:val="auto"/></w:rPr><w:t>cin</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b w:val="0"/><w:color w:val="auto"/></w:rPr><w:t>)</w:t></w:r>'

Replace-RunXml $d " that allows your CPU to read text/numbers from a keyboard buffer(kind of like cin)" $cinNew | Out-Null

# ---------------------------------------------------------------------
# Edit 11: footer page-number field cached result "2" -> "1"
# ---------------------------------------------------------------------

$footer = $d.Sections(1).Footers(1)
$footer.Range.Characters(1).Text = "1"

Write-Host "all edits applied"
